# This workbook is a periodic COVID-19 "paises" (countries) data dump, sorted
# descending by "Casos totales" (column B). This edit refreshes the snapshot:
# the "Datos actualizados" timestamp moves from 17:03 to 18:03, many countries
# get updated totals, and a handful of countries change relative rank (so their
# row keeps its position in the sort order but shows a different country).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 6 de Mayo de 2020 a las 18:03"

# row -> updated values for columns A (country) through H (muertes)
$rowUpdates = @{
    4 = @("Estados Unidos", 1241430, 3797, 201879, 966835, 16179, 445, 72716)
    6 = @("Italia", 214457, 1444, 93245, 91528, 1333, 369, 29684)
    16 = @("India", 52340, 2940, 14911, 35661, 0, 75, 1768)
    17 = @("Peru", 51189, 0, 15413, 34332, 709, 0, 1444)
    18 = @("Belgica", 50781, 272, 12731, 29711, 646, 323, 8339)
    29 = @("Singapur", 20198, 788, 1634, 18544, 23, 2, 20)
    36 = @("Polonia", 14740, 309, 4655, 9352, 160, 17, 733)
    48 = @("Chequia", 7933, 37, 4202, 3469, 59, 5, 262)
    49 = @("Egipto", 7588, 387, 1815, 5304, 0, 17, 469)
    50 = @("Sudafrica", 7572, 0, 2746, 4678, 36, 0, 148)
    51 = @("Panama", 7523, 136, 823, 6490, 88, 7, 210)
    58 = @("Argelia", 4997, 159, 2197, 2324, 22, 6, 476)
    70 = @("Grecia", 2663, 21, 1374, 1142, 36, 1, 147)
    82 = @("Cuba", 1703, 18, 1001, 633, 8, 0, 69)
    95 = @("Republica de Chipre", 883, 5, 296, 572, 15, 0, 15)
    96 = @("Somalia", 873, 38, 87, 747, 2, 1, 39)
    97 = @("Kirguistan", 871, 28, 614, 245, 13, 1, 12)
    98 = @("Sudan", 852, 74, 80, 727, 0, 0, 45)
    101 = @("Sri Lanka", 795, 24, 215, 571, 1, 0, 9)
    112 = @("Mali", 631, 19, 261, 338, 0, 0, 32)
    113 = @("Maldivas", 618, 45, 20, 596, 2, 0, 2)
    116 = @("Kenia", 582, 47, 190, 366, 1, 2, 26)
    121 = @("Jordania", 473, 2, 377, 87, 5, 0, 9)
    140 = @("Liberia", 178, 8, 75, 83, 0, 0, 20)
    141 = @("Santo Tome y Principe", 174, 0, 4, 167, 0, 0, 3)
    142 = @("Republica del Chad", 170, 0, 43, 110, 0, 0, 17)
    175 = @("Malaui", 43, 2, 9, 31, 1, 0, 3)
    198 = @("Dominica", 16, 0, 14, 2, 0, 0, 0)
    199 = @("Curazao", 16, 0, 13, 2, 0, 0, 1)
    205 = @("Seychelles", 11, 0, 8, 3, 0, 0, 0)
    206 = @("Montserrat", 11, 0, 7, 3, 1, 0, 1)
}

foreach ($row in $rowUpdates.Keys) {
    $values = $rowUpdates[$row]
    for ($col = 0; $col -lt $values.Length; $col++) {
        $ws.Cells.Item([int]$row, $col + 1).Value = $values[$col]
    }
}
